# Insert a new weekly price record as row 47 in the "Haba" sheet,
# pushing the former rows 47-54 down to 48-55.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 47 (shifts existing rows 47-54 -> 48-55)
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with the new record's data
$ws.Cells.Item(47, 1).Value  = 10
$ws.Cells.Item(47, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value  = "La Araucanía"
$ws.Cells.Item(47, 4).Value  = 44522
$ws.Cells.Item(47, 5).Value  = 9
$ws.Cells.Item(47, 6).Value  = 100112026
$ws.Cells.Item(47, 7).Value  = "Haba"
$ws.Cells.Item(47, 8).Value  = "Sin especificar"
$ws.Cells.Item(47, 9).Value  = "Primera"
$ws.Cells.Item(47, 10).Value = 40
$ws.Cells.Item(47, 11).Value = 8000
$ws.Cells.Item(47, 12).Value = 8000
$ws.Cells.Item(47, 13).Value = 8000
$ws.Cells.Item(47, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(47, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(47, 16).Value = 320
$ws.Cells.Item(47, 17).Value = 25
$ws.Cells.Item(47, 18).Value = "Hortaliza"
